# Revert capacity charts to show kilowatts on the y-axis.
# - Divide the "Waste Gas" column (E) data for rows 15-26 by 1000 (Watts -> Kilowatts)
# - Update the number format used by those cells to show one decimal place
# - Update the chart's value-axis title and number format accordingly

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New kilowatt values (Watts / 1000) for Sheet1!E15:E26
$newValues = @{
    15 = 4.1
    16 = 9
    17 = 31
    18 = 13.8
    19 = 0
    20 = 41.4
    21 = 63.1
    22 = 49
    23 = 124.233
    24 = 141.43
    25 = 27.96
    26 = 33.932
}

foreach ($row in $newValues.Keys) {
    $ws.Range("E$row").Value = $newValues[$row]
}

# Update the shared number format (numFmtId 164, applied to these cells) to one
# decimal place so fractional kilowatt values are visible.
$ws.Range("E15:E26").NumberFormat = "#,##0.0"

# Update the chart: axis title text and value-axis number format.
$chartObj = $ws.ChartObjects().Item(1)
$chart = $chartObj.Chart
$valueAxis = $chart.Axes(2)  # xlValue
$valueAxis.AxisTitle.Text = "Kilowatts (kW)"
$valueAxis.TickLabels.NumberFormat = "#,##0"
